$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new displayed value. Values are written with a leading apostrophe
# (text-prefix) so Excel stores them as literal text instead of re-parsing
# numeric-looking strings ("330.24", "1.31%", ...) into Number/Percentage
# cells, which would silently drop significant trailing zeros (e.g. the
# "0.0001200" / "0.00000000750" rows).
$changes = @{
    'D2' = '330.24'
    'E2' = '1.31%'
    'D3' = '44.29'
    'E3' = '0.29%'
    'D4' = '5.493'
    'E4' = '-1.69%'
    'D5' = '0.08038'
    'E5' = '-0.08%'
    'D6' = '2.017'
    'E6' = '6.95%'
    'D7' = '0.9527'
    'E7' = '0.85%'
    'D8' = '2.560'
    'E8' = '-3.29%'
    'D9' = '0.1145'
    'E9' = '-1.40%'
    'D10' = '0.1899'
    'E10' = '3.39%'
    'D11' = '10.75'
    'E11' = '28.34%'
    'D12' = '0.09832'
    'E12' = '-0.17%'
    'D13' = '0.04825'
    'E13' = '11.86%'
    'D14' = '0.1065'
    'E14' = '-0.08%'
    'D15' = '0.001281'
    'E15' = '0.57%'
    'D16' = '0.04073'
    'D17' = '0.005880'
    'E17' = '-0.98%'
    'D18' = '3.367'
    'E18' = '-6.68%'
    'D19' = '4.402'
    'E19' = '2.41%'
    'D20' = '0.3431'
    'E20' = '-1.86%'
    'E21' = '1.44%'
    'E22' = '-5.75%'
    'D23' = '0.001273'
    'E23' = '2.18%'
    'D24' = '0.004361'
    'E24' = '-3.14%'
    'D25' = '0.0001200'
    'E25' = '-4.96%'
    'D26' = '0.0003741'
    'E26' = '-6.37%'
    'D38' = '0.02603'
    'E38' = '-0.76%'
    'D39' = '0.05774'
    'E39' = '5.81%'
    'D40' = '0.007551'
    'E40' = '-0.88%'
    'E41' = '0.77%'
    'D42' = '0.007143'
    'E42' = '-2.59%'
    'D43' = '0.002014'
    'E43' = '-0.25%'
    'D44' = '0.008828'
    'E44' = '0.10%'
    'D45' = '0.00006977'
    'E45' = '0.97%'
    'D46' = '0.00000000750'
    'E46' = '-0.20%'
    'E47' = '-0.28%'
    'E48' = '55.12%'
    'D49' = '0.003552'
    'E49' = '-3.48%'
    'D50' = '0.00002099'
    'E50' = '-0.20%'
    'D51' = '0.0001999'
    'E51' = '-0.20%'
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = "'" + $changes[$addr]
}
